$wb = $excel.ActiveWorkbook

# --- Duplicate "Sheet2" so we get both the renamed "July 8th to 21st" sheet
#     (new data) and a "Sheet2 (2)" copy that keeps the old June 10-23 blank
#     template data, exactly like Excel's own "Move or Copy... (Create a
#     copy)" workflow.
$ws4 = $wb.Worksheets.Item("Sheet2")
$ws4.Unprotect()
$ws4.Copy([System.Reflection.Missing]::Value, $ws4)

$wsCopy = $wb.Worksheets.Item("Sheet2 (2)")

# Rename the original tab and point it at the next two-week period.
$ws4.Name = "July 8th to 21st"

# --- Update the dates for July 8th to 21st (formerly June 10-23).
$dates = @(43654,43655,43656,43657,43658,43659,43660,43661,43662,43663,43664,43665,43666,43667)
for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws4.Cells.Item($r, 2).Value = $dates[$i]
}

# --- Fill in Monday's (7/8) IN/OUT time punch.
$ws4.Range("C2:D2").NumberFormat = "h:mm AM/PM"
$ws4.Range("C2").Value = 0.4375
$ws4.Range("D2").Value = 0.72916666666666663

# --- The worksheet copy doesn't bring its ListObject/Table along, so
#     re-create it on "Sheet2 (2)" with the same range/columns/style Excel
#     would have produced for the duplicated tab.
$lo = $wsCopy.ListObjects.Add(1, $wsCopy.Range("A1:G15"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table146"
$lo.TableStyle = "TableStyleLight1"

# --- Selection / active-tab bookkeeping so it matches the recorded UI state.
$wsCopy.Activate()
$wsCopy.Range("D15").Select()

$ws4.Activate()
$ws4.Range("C3").Select()
